$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2261.1667
$ws.Range("J17").Value = 2261.1667
$ws.Range("L17").Value = 6783.500100000001
$ws.Range("N17").Value = -7119.500100000001

$ws.Range("H19").Value = 195.13333
$ws.Range("I19").Value = 181.77777
$ws.Range("J19").Value = 215.16667
$ws.Range("K19").Value = 181.77777
$ws.Range("L19").Value = 215.16667
$ws.Range("M19").Value = -6.777770000000004
$ws.Range("N19").Value = -565.1666700000001

$ws.Range("H106").Value = 2711.375
$ws.Range("I106").Value = 1982.1666
$ws.Range("K106").Value = 1982.1666
$ws.Range("M106").Value = -1351.1666

$ws.Range("H132").Value = 1201
$ws.Range("I132").Value = 1264.2808
$ws.Range("J132").Value = 299.25
$ws.Range("K132").Value = 3792.8424
$ws.Range("L132").Value = 897.75
$ws.Range("M132").Value = -1262.8424
$ws.Range("N132").Value = -5957.75

$ws.Range("H135").Value = 1249.7941
$ws.Range("I135").Value = 1249.7941
$ws.Range("K135").Value = 11248.1469
$ws.Range("M135").Value = -8713.1469

$ws.Range("H137").Value = 1244311.6
$ws.Range("I137").Value = 1389711.9
$ws.Range("J137").Value = 1194173.6
$ws.Range("K137").Value = 4169135.7
$ws.Range("L137").Value = 3582520.8
$ws.Range("M137").Value = -4166585.7
$ws.Range("N137").Value = -3587620.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 2547524
$ws.Range("I8").Value = 2801026.5
$ws.Range("K8").Value = 2801026.5
$ws.Range("M8").Value = -2800882.5

$ws.Range("H13").Value = 6057.25
$ws.Range("I13").Value = 615
$ws.Range("J13").Value = 11499.5
$ws.Range("K13").Value = 615
$ws.Range("L13").Value = 11499.5
$ws.Range("M13").Value = -471
$ws.Range("N13").Value = -11787.5

$ws.Range("H32").Value = 2832.9275
$ws.Range("I32").Value = 2801.0588
$ws.Range("K32").Value = 2801.0588
$ws.Range("M32").Value = -2514.0588

$ws.Range("H45").Value = 3260.36
$ws.Range("I45").Value = 3195.5
$ws.Range("K45").Value = 3195.5
$ws.Range("M45").Value = -2818.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1879.1305
$ws.Range("I20").Value = 2111.7273
$ws.Range("K20").Value = 2111.7273
$ws.Range("M20").Value = -1864.7273

$ws.Range("H82").Value = 51149.273
$ws.Range("I82").Value = 16731.666
$ws.Range("J82").Value = 92450.39999999999
$ws.Range("K82").Value = 16731.666
$ws.Range("L82").Value = 92450.39999999999
$ws.Range("M82").Value = -16348.666
$ws.Range("N82").Value = -93216.39999999999

$ws.Range("H85").Value = 51149.273
$ws.Range("I85").Value = 16731.666
$ws.Range("J85").Value = 92450.39999999999
$ws.Range("K85").Value = 16731.666
$ws.Range("L85").Value = 92450.39999999999
$ws.Range("M85").Value = -15405.666
$ws.Range("N85").Value = -95102.39999999999

$ws.Range("H134").Value = 20696.492
$ws.Range("I134").Value = 23968.455
$ws.Range("J134").Value = 9622.154
$ws.Range("K134").Value = 71905.36500000001
$ws.Range("L134").Value = 28866.462
$ws.Range("M134").Value = -69370.36500000001
$ws.Range("N134").Value = -33936.462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 166963.94
$ws.Range("I7").Value = 250159.5
$ws.Range("K7").Value = 250159.5
$ws.Range("M7").Value = -250046.5

$ws.Range("H22").Value = 1081.091
$ws.Range("I22").Value = 999.3
$ws.Range("K22").Value = 999.3
$ws.Range("M22").Value = -649.3

$ws.Range("H99").Value = 8256.733
$ws.Range("I99").Value = 4455.2
$ws.Range("J99").Value = 12058.267
$ws.Range("K99").Value = 4455.2
$ws.Range("L99").Value = 12058.267
$ws.Range("M99").Value = -2957.2
$ws.Range("N99").Value = -15054.267

$ws.Range("H107").Value = 1948
$ws.Range("I107").Value = 1337.3334
$ws.Range("K107").Value = 1337.3334
$ws.Range("M107").Value = 582.6666

$ws.Range("H119").Value = 98408.8
$ws.Range("J119").Value = 98408.8
$ws.Range("L119").Value = 98408.8
$ws.Range("N119").Value = -108084.8

$ws.Range("H126").Value = 8256.733
$ws.Range("I126").Value = 4455.2
$ws.Range("J126").Value = 12058.267
$ws.Range("K126").Value = 13365.6
$ws.Range("L126").Value = 36174.801
$ws.Range("M126").Value = -10895.6
$ws.Range("N126").Value = -41114.801

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 118189.06
$ws.Range("I14").Value = 118189.06
$ws.Range("K14").Value = 354567.18
$ws.Range("M14").Value = -354394.18

$ws.Range("H129").Value = 4858
$ws.Range("I129").Value = 1699.75
$ws.Range("K129").Value = 5099.25
$ws.Range("M129").Value = -99.25

$ws.Range("H131").Value = 1673.3478
$ws.Range("I131").Value = 799.7143
$ws.Range("J131").Value = 2055.5625
$ws.Range("K131").Value = 2399.1429
$ws.Range("L131").Value = 6166.6875
$ws.Range("M131").Value = 2640.8571
$ws.Range("N131").Value = -16246.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 21500
$ws.Range("I48").Value = 21500
$ws.Range("K48").Value = 21500
$ws.Range("M48").Value = -21015

$ws.Range("H102").Value = 4569.4414
$ws.Range("I102").Value = 4512.033
$ws.Range("K102").Value = 4512.033
$ws.Range("M102").Value = -2890.033

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2466.3333
$ws.Range("I7").Value = 2466.3333
$ws.Range("K7").Value = 2466.3333
$ws.Range("M7").Value = -2354.3333

$ws.Range("H12").Value = 2036.1818
$ws.Range("I12").Value = 1928.5714
$ws.Range("J12").Value = 2224.5
$ws.Range("K12").Value = 1928.5714
$ws.Range("L12").Value = 2224.5
$ws.Range("M12").Value = -1758.5714
$ws.Range("N12").Value = -2564.5

$ws.Range("H22").Value = 55416
$ws.Range("I22").Value = 127317.5
$ws.Range("J22").Value = 3124
$ws.Range("K22").Value = 127317.5
$ws.Range("L22").Value = 3124
$ws.Range("M22").Value = -127022.5
$ws.Range("N22").Value = -3714

$ws.Range("H27").Value = 55416
$ws.Range("I27").Value = 127317.5
$ws.Range("J27").Value = 3124
$ws.Range("K27").Value = 127317.5
$ws.Range("L27").Value = 3124
$ws.Range("M27").Value = -127210.5
$ws.Range("N27").Value = -3338

$ws.Range("H40").Value = 7659.4443
$ws.Range("I40").Value = 6276.4287
$ws.Range("K40").Value = 6276.4287
$ws.Range("M40").Value = -6140.4287

$ws.Range("H46").Value = 939.85
$ws.Range("I46").Value = 980.0714
$ws.Range("K46").Value = 980.0714
$ws.Range("M46").Value = -792.0714

$ws.Range("H61").Value = 3481.45
$ws.Range("I61").Value = 2018.75
$ws.Range("J61").Value = 5675.5
$ws.Range("K61").Value = 2018.75
$ws.Range("L61").Value = 5675.5
$ws.Range("M61").Value = -1816.75
$ws.Range("N61").Value = -6079.5

$ws.Range("H113").Value = 3481.45
$ws.Range("I113").Value = 2018.75
$ws.Range("J113").Value = 5675.5
$ws.Range("K113").Value = 2018.75
$ws.Range("L113").Value = 5675.5
$ws.Range("M113").Value = 151.25
$ws.Range("N113").Value = -10015.5

$ws.Range("H122").Value = 3214.9285
$ws.Range("I122").Value = 2957.1292
$ws.Range("J122").Value = 3941.4546
$ws.Range("K122").Value = 8871.3876
$ws.Range("L122").Value = 11824.3638
$ws.Range("M122").Value = -6421.3876
$ws.Range("N122").Value = -16724.3638

$ws.Range("H126").Value = 2466.3333
$ws.Range("I126").Value = 2466.3333
$ws.Range("K126").Value = 7398.999899999999
$ws.Range("M126").Value = -4928.999899999999

$ws.Range("H130").Value = 68970.82000000001
$ws.Range("J130").Value = 68970.82000000001
$ws.Range("L130").Value = 68970.82000000001
$ws.Range("N130").Value = -79010.82000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 34500
$ws.Range("I3").Value = 34000
$ws.Range("J3").Value = 35000
$ws.Range("K3").Value = 34000
$ws.Range("L3").Value = 35000
$ws.Range("M3").Value = -33886
$ws.Range("N3").Value = -35228

$ws.Range("H13").Value = 2226.25
$ws.Range("I13").Value = 1001.6667
$ws.Range("K13").Value = 1001.6667
$ws.Range("M13").Value = -861.6667

$ws.Range("H126").Value = 2460.4
$ws.Range("I126").Value = 1833.75
$ws.Range("K126").Value = 5501.25
$ws.Range("M126").Value = -3031.25
